$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 with new cluster counts (Birch + K-means method)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 243

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 96

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = 96

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = 10

# Remove row 6 entirely (it no longer exists in the updated data)
$ws.Range("A6:B6").Delete()
